$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate row 8 through row 15 with inspection plan data
$ws.Range("C8").Value = 'MLK_PMT_10103_-_V-003'
$ws.Range("D8").Value = 'Air Receiver'
$ws.Range("E8").Value = 'Bottom Head'
$ws.Range("G8").Value = 'CONDENSATE'
$ws.Range("H8").Value = 'Not Found'
$ws.Range("I8").Value = 'A/SA 516'
$ws.Range("J8").Value = 'Gr.70'
$ws.Range("K8").Value = 'N/A'
$ws.Range("L8").Value = '200°C'
$ws.Range("M8").Value = '1,000 kPaG'
$ws.Range("N8").Value = '185°C'
$ws.Range("O8").Value = '1,000 kPaG'

$ws.Range("E9").Value = 'Top Head'
$ws.Range("G9").Value = 'CONDENSATE'
$ws.Range("H9").Value = 'Not Found'
$ws.Range("I9").Value = 'A/SA 105'
$ws.Range("K9").Value = 'N/A'
$ws.Range("L9").Value = '200°C'
$ws.Range("M9").Value = '1,000 kPaG'
$ws.Range("N9").Value = '185°C'
$ws.Range("O9").Value = '1,000 kPaG'

$ws.Range("E10").Value = 'Shell'
$ws.Range("G10").Value = 'CONDENSATE'
$ws.Range("H10").Value = 'Not Found'
$ws.Range("I10").Value = 'A/SA 106'
$ws.Range("J10").Value = 'Gr.B'
$ws.Range("K10").Value = 'N/A'
$ws.Range("L10").Value = '200°C'
$ws.Range("M10").Value = '1,000 kPaG'
$ws.Range("N10").Value = '185°C'
$ws.Range("O10").Value = '1,000 kPaG'

$ws.Range("E11").Value = 'Socket'
$ws.Range("G11").Value = 'CONDENSATE'
$ws.Range("H11").Value = 'Not Found'
$ws.Range("I11").Value = 'A/SA 105'
$ws.Range("K11").Value = 'N/A'
$ws.Range("L11").Value = '200°C'
$ws.Range("M11").Value = '1,000 kPaG'
$ws.Range("N11").Value = '185°C'
$ws.Range("O11").Value = '1,000 kPaG'

$ws.Range("E12").Value = 'Tube Bundle'
$ws.Range("G12").Value = 'CONDENSATE'
$ws.Range("H12").Value = 'Not Found'
$ws.Range("I12").Value = 'A/SA 516'
$ws.Range("J12").Value = 'Gr.70'
$ws.Range("K12").Value = 'N/A'
$ws.Range("L12").Value = '200°C'
$ws.Range("M12").Value = '1,000 kPaG'
$ws.Range("N12").Value = '185°C'
$ws.Range("O12").Value = '1,000 kPaG'

$ws.Range("E13").Value = 'Head'
$ws.Range("G13").Value = 'CONDENSATE'
$ws.Range("H13").Value = 'Structural Steel'
$ws.Range("I13").Value = 'S275JR'
$ws.Range("K13").Value = 'N/A'
$ws.Range("L13").Value = '200°C'
$ws.Range("M13").Value = '1,000 kPaG'
$ws.Range("N13").Value = '185°C'
$ws.Range("O13").Value = '1,000 kPaG'

$ws.Range("E14").Value = 'Pressure Retaining Bolt & Nut'
$ws.Range("G14").Value = 'CONDENSATE'
$ws.Range("H14").Value = 'Stainless Steel Bolting'
$ws.Range("I14").Value = 'SA193 / SA194'
$ws.Range("J14").Value = 'B7 / 2H'
$ws.Range("K14").Value = 'N/A'
$ws.Range("L14").Value = '200°C'
$ws.Range("M14").Value = '1,000 kPaG'
$ws.Range("N14").Value = '185°C'
$ws.Range("O14").Value = '1,000 kPaG'

$ws.Range("E15").Value = 'External Fittings'
$ws.Range("G15").Value = 'CONDENSATE'
$ws.Range("H15").Value = 'Carbon Steel'
$ws.Range("I15").Value = 'JIS G3507'
$ws.Range("K15").Value = 'N/A'
$ws.Range("L15").Value = '200°C'
$ws.Range("M15").Value = '1,000 kPaG'
$ws.Range("N15").Value = '185°C'
$ws.Range("O15").Value = '1,000 kPaG'

# Fix merged cell ranges for columns A, B, C, D (rows 8-38 -> rows 8-15)
$ws.Range("A8:A38").UnMerge()
$ws.Range("B8:B38").UnMerge()
$ws.Range("C8:C38").UnMerge()
$ws.Range("D8:D38").UnMerge()
$ws.Range("A8:A15").Merge()
$ws.Range("B8:B15").Merge()
$ws.Range("C8:C15").Merge()
$ws.Range("D8:D15").Merge()
